# TC for Manage Exclusion Date
# Insert a new test-case row for "Manage Exclusion Date" into the
# "Test Result" sheet, right after the "Manage Guardian UI Verification"
# section header (i.e. before the old row 6), shifting every row below
# it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Result")

# Insert a new blank row at row 6 (everything from row 6 down shifts to 7+)
$ws.Rows("6:6").Insert()

# Copy the formatting of the row above (row 5) into the new row 6 so the
# borders/styles match the surrounding "FUNCTION" rows (A:s1 / B:s27 / C:s30 / D:s26)
$ws.Range("A5:D5").Copy()
$ws.Range("A6:D6").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Match row height / thick-bottom-border sizing used by the neighbouring rows
$ws.Rows("6:6").RowHeight = $ws.Rows("5:5").RowHeight

# Fill in the new scenario row
$ws.Range("B6").Value = "Manage Exclusion Date"
$ws.Range("C6").Value = $false

# The two previously-checked ("Print Labels" / "Receive") checkboxes are
# unchecked again for this new snapshot of results (now rows 19 & 20 after
# the insertion above).
$ws.Range("C19").Value = $false
$ws.Range("C20").Value = $false

# Update the conditional-formatting range so it still covers the whole
# results column (C2:C42 -> C2:C43)
$fc = $ws.Range("C2:C42").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("C2:C43"))

# Update selection / view to match the saved state referenced by the diff
$ws.Range("D15").Select()
